$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the 'Desenho tecnico mecanico' class entries from segunda/terça (B2, C4)
# to sexta (F3, F4), per the requested update ("Só para o Túlio pegar atualizado").

$ws.Range("B2").Value = "-"
$ws.Range("F3").Value = "['MEC-1B-Desenho tecnico mecanico', 'MEC-1B-Desenho tecnico mecanico']"
$ws.Range("C4").Value = "-"
$ws.Range("F4").Value = "[-, 'MEC-1B-Desenho tecnico mecanico']"
